$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting existing rows (2-5) down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new "machine" header names.
$ws.Range("A2").Value = "orden"
$ws.Range("B2").Value = "comarca-nombre"
$ws.Range("C2").Value = "comarca-codigo"
$ws.Range("D2").Value = "siglas"
$ws.Range("E2").Value = "diputados"
$ws.Range("F2").Value = "provincia-codigo"
$ws.Range("G2").Value = "municipio-codigo"
$ws.Range("H2").Value = "provincia-nombre"
$ws.Range("I2").Value = "ano"
$ws.Range("J2").Value = "votos"
$ws.Range("K2").Value = "municipio-nombre"

# Remove the old trailing row that only held "mapping-ano.xlsx" in I (now row 6 after insert).
$ws.Rows.Item(6).Delete()
